$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 3) mirrors the structure of row 2
$ws.Range("A3").Value = "test@gmail.com"
$ws.Range("B3").Value = "Ram"
$ws.Range("C3").Value = "Rddy"
$ws.Range("D3").Value = "Myadd2"
$ws.Range("E3").Value = "Amaravathi"
$ws.Range("F3").Value = "IN"
$ws.Range("G3").Value = "ap"
$ws.Range("H3").Value = 123456
$ws.Range("I3").Value = 1234567899
$ws.Range("J3").Value = $true

# Add hyperlink on the email cell, matching A2's hyperlink setup
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:test@gmail.com")
$ws.Range("A3").Style = "Hyperlink"

$ws.Range("J3").Select()
